# "fixed excel data from SQL": the workbook originally stopped at row 201.
# This adds 55 more country/aggregate rows (pulled from the same SQL source,
# already sorted desc. by percentpopulationinfected) below the existing table,
# re-applies the AutoFilter + its companion hidden _FilterDatabase defined name
# over the original A1:D201 extent, and leaves the active cell where the author
# left it (F9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply AutoFilter (and the hidden _FilterDatabase name it writes) over the
# original table extent while it is still the full used range, so the filter
# ref is not later auto-expanded once the extra rows below are populated.
$null = $ws.Range("A1:D201").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$201")
$filterName.Visible = $false

# Column A (location) is written cell-by-cell in the exact order the values were
# originally entered/imported (so the shared-string table ends up in the same
# order as the source workbook) even though the rows below were later sorted by
# percentpopulationinfected, which put row 210 (Africa) ahead of rows 211-212.
$ws.Range("A202").Value = "Lesotho"
$ws.Range("A203").Value = "Djibouti"
$ws.Range("A204").Value = "Mauritania"
$ws.Range("A205").Value = "Bangladesh"
$ws.Range("A206").Value = "Myanmar"
$ws.Range("A207").Value = "Comoros"
$ws.Range("A208").Value = "Equatorial Guinea"
$ws.Range("A209").Value = "Rwanda"
$ws.Range("A211").Value = "Cambodia"
$ws.Range("A212").Value = "Uzbekistan"
$ws.Range("A210").Value = "Africa"
$ws.Range("A213").Value = "Mozambique"
$ws.Range("A214").Value = "Pakistan"
$ws.Range("A215").Value = "Kenya"
$ws.Range("A216").Value = "Algeria"
$ws.Range("A217").Value = "Senegal"
$ws.Range("A218").Value = "Ghana"
$ws.Range("A219").Value = "Afghanistan"
$ws.Range("A220").Value = "Gambia"
$ws.Range("A221").Value = "Egypt"
$ws.Range("A222").Value = "Papua New Guinea"
$ws.Range("A223").Value = "Cameroon"
$ws.Range("A224").Value = "Togo"
$ws.Range("A225").Value = "Malawi"
$ws.Range("A226").Value = "Guinea-Bissau"
$ws.Range("A227").Value = "Congo"
$ws.Range("A228").Value = "Burundi"
$ws.Range("A229").Value = "Ethiopia"
$ws.Range("A230").Value = "Uganda"
$ws.Range("A231").Value = "Cote d'Ivoire"
$ws.Range("A232").Value = "Low income"
$ws.Range("A233").Value = "Angola"
$ws.Range("A234").Value = "Haiti"
$ws.Range("A235").Value = "Eritrea"
$ws.Range("A236").Value = "Guinea"
$ws.Range("A237").Value = "Central African Republic"
$ws.Range("A238").Value = "Tokelau"
$ws.Range("A239").Value = "Syria"
$ws.Range("A240").Value = "Madagascar"
$ws.Range("A241").Value = "Nicaragua"
$ws.Range("A242").Value = "Benin"
$ws.Range("A243").Value = "Tajikistan"
$ws.Range("A244").Value = "South Sudan"
$ws.Range("A245").Value = "Somalia"
$ws.Range("A246").Value = "Liberia"
$ws.Range("A247").Value = "Mali"
$ws.Range("A248").Value = "Sudan"
$ws.Range("A249").Value = "Nigeria"
$ws.Range("A250").Value = "Burkina Faso"
$ws.Range("A251").Value = "Democratic Republic of Congo"
$ws.Range("A252").Value = "Sierra Leone"
$ws.Range("A253").Value = "Tanzania"
$ws.Range("A254").Value = "Chad"
$ws.Range("A255").Value = "Niger"
$ws.Range("A256").Value = "Yemen"

# Columns B:D (population, highestinfectioncount, percentpopulationinfected) for
# the new rows 202-256, in row order.
$data = New-Object 'object[,]' 55,3
$data[0,0] = 2305826
$data[0,1] = 34490
$data[0,2] = 1.4957763509
$data[1,0] = 1120851
$data[1,1] = 15690
$data[1,2] = 1.3998292369
$data[2,0] = 4736146
$data[2,1] = 63439
$data[2,2] = 1.3394646195
$data[3,0] = 171186368
$data[3,1] = 2037947
$data[3,2] = 1.1904843965
$data[4,0] = 54179312
$data[4,1] = 633967
$data[4,2] = 1.1701274464
$data[5,0] = 836783
$data[5,1] = 9048
$data[5,2] = 1.0812839171
$data[6,0] = 1674916
$data[6,1] = 17130
$data[6,2] = 1.0227378567
$data[7,0] = 13776702
$data[7,1] = 133194
$data[7,2] = 0.966806134
$data[8,0] = 1426736614
$data[8,1] = 13063184
$data[8,2] = 0.9155988479
$data[9,0] = 16767851
$data[9,1] = 138720
$data[9,2] = 0.8272974277
$data[10,0] = 34627648
$data[10,1] = 251430
$data[10,2] = 0.7260960952
$data[11,0] = 32969520
$data[11,1] = 233214
$data[11,2] = 0.7073624366
$data[12,0] = 235824864
$data[12,1] = 1577699
$data[12,2] = 0.6690130011
$data[13,0] = 54027484
$data[13,1] = 342943
$data[13,2] = 0.634756562
$data[14,0] = 44903228
$data[14,1] = 271522
$data[14,2] = 0.6046825854
$data[15,0] = 17316452
$data[15,1] = 88933
$data[15,2] = 0.5135751827
$data[16,0] = 33475870
$data[16,1] = 171281
$data[16,2] = 0.5116551116
$data[17,0] = 41128772
$data[17,1] = 209602
$data[17,2] = 0.5096237738
$data[18,0] = 2705995
$data[18,1] = 12598
$data[18,2] = 0.4655588795
$data[19,0] = 110990096
$data[19,1] = 515792
$data[19,2] = 0.4647189421
$data[20,0] = 10142625
$data[20,1] = 46826
$data[20,2] = 0.4616753552
$data[21,0] = 27914542
$data[21,1] = 124605
$data[21,2] = 0.4463802415
$data[22,0] = 8848700
$data[22,1] = 39407
$data[22,2] = 0.4453422537
$data[23,0] = 20405318
$data[23,1] = 88613
$data[23,2] = 0.4342642443
$data[24,0] = 2105580
$data[24,1] = 8960
$data[24,2] = 0.4255359568
$data[25,0] = 5970430
$data[25,1] = 25110
$data[25,2] = 0.4205727226
$data[26,0] = 12889583
$data[26,1] = 53661
$data[26,2] = 0.4163129249
$data[27,0] = 123379928
$data[27,1] = 500169
$data[27,2] = 0.4053892786
$data[28,0] = 47249588
$data[28,1] = 170463
$data[28,2] = 0.3607713997
$data[29,0] = 28160548
$data[29,1] = 88277
$data[29,2] = 0.3134775644
$data[30,0] = 737604900
$data[30,1] = 2287803
$data[30,2] = 0.310166459
$data[31,0] = 35588996
$data[31,1] = 105298
$data[31,2] = 0.2958723534
$data[32,0] = 11585003
$data[32,1] = 34202
$data[32,2] = 0.2952265097
$data[33,0] = 3684041
$data[33,1] = 10189
$data[33,2] = 0.2765712977
$data[34,0] = 13859349
$data[34,1] = 38280
$data[34,2] = 0.2762034494
$data[35,0] = 5579148
$data[35,1] = 15367
$data[35,2] = 0.2754363211
$data[36,0] = 1893
$data[36,1] = 5
$data[36,2] = 0.264131009
$data[37,0] = 22125242
$data[37,1] = 57423
$data[37,2] = 0.2595361443
$data[38,0] = 29611718
$data[38,1] = 67941
$data[38,2] = 0.2294395752
$data[39,0] = 6948395
$data[39,1] = 15672
$data[39,2] = 0.22554849
$data[40,0] = 13352864
$data[40,1] = 27999
$data[40,2] = 0.2096853529
$data[41,0] = 9952789
$data[41,1] = 17786
$data[41,2] = 0.1787036779
$data[42,0] = 10913172
$data[42,1] = 18368
$data[42,2] = 0.1683103684
$data[43,0] = 17597508
$data[43,1] = 27324
$data[43,2] = 0.1552719851
$data[44,0] = 5302690
$data[44,1] = 8090
$data[44,2] = 0.152564076
$data[45,0] = 22593598
$data[45,1] = 33067
$data[45,2] = 0.1463556181
$data[46,0] = 46874200
$data[46,1] = 63853
$data[46,2] = 0.1362220582
$data[47,0] = 218541216
$data[47,1] = 266641
$data[47,2] = 0.1220094794
$data[48,0] = 22673764
$data[48,1] = 22056
$data[48,2] = 0.0972754237
$data[49,0] = 99010216
$data[49,1] = 95814
$data[49,2] = 0.0967718321
$data[50,0] = 8605723
$data[50,1] = 7760
$data[50,2] = 0.0901725515
$data[51,0] = 65497752
$data[51,1] = 42927
$data[51,2] = 0.0655396539
$data[52,0] = 17723312
$data[52,1] = 7682
$data[52,2] = 0.0433440431
$data[53,0] = 26207982
$data[53,1] = 9513
$data[53,2] = 0.0362981019
$data[54,0] = 33696612
$data[54,1] = 11945
$data[54,2] = 0.0354486677
$ws.Range("B202:D256").Value = $data

# Restore the active selection to F9.
$null = $ws.Range("F9").Select()
